# Auto-generated Excel COM-interop script
# Applies numeric updates to the "Sheets" (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR) per scheduled-runner refresh.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2223.3845
$ws.Range("I40").Value = 1900.0
$ws.Range("J40").Value = 2600.6667
$ws.Range("K40").Value = 1900.0
$ws.Range("L40").Value = 2600.6667
$ws.Range("M40").Value = -1725.0
$ws.Range("N40").Value = -2950.6667
$ws.Range("H112").Value = 1027.0952
$ws.Range("J112").Value = 1053.95
$ws.Range("L112").Value = 3161.85
$ws.Range("N112").Value = -5377.85
$ws.Range("H113").Value = 168800.83
$ws.Range("I113").Value = 252251.25
$ws.Range("J113").Value = 1900.0
$ws.Range("K113").Value = 252251.25
$ws.Range("L113").Value = 1900.0
$ws.Range("M113").Value = -248997.25
$ws.Range("N113").Value = -8408.0
$ws.Range("H129").Value = 2314.2922
$ws.Range("I129").Value = 5470.75
$ws.Range("J129").Value = 911.42224
$ws.Range("K129").Value = 16412.25
$ws.Range("L129").Value = 2734.26672
$ws.Range("M129").Value = -11412.25
$ws.Range("N129").Value = -12734.26672
$ws.Range("H132").Value = 4812560.5
$ws.Range("I132").Value = 5686908.0
$ws.Range("J132").Value = 3650.375
$ws.Range("K132").Value = 17060724.0
$ws.Range("L132").Value = 10951.125
$ws.Range("M132").Value = -17058194.0
$ws.Range("N132").Value = -16011.125

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H46").Value = 3298.8
$ws.Range("H58").Value = 12759.8
$ws.Range("J58").Value = 12759.8
$ws.Range("L58").Value = 12759.8
$ws.Range("N58").Value = -13619.8
$ws.Range("H74").Value = 845.96295
$ws.Range("I74").Value = 741.7222
$ws.Range("J74").Value = 1054.4445
$ws.Range("K74").Value = 741.7222
$ws.Range("L74").Value = 1054.4445
$ws.Range("M74").Value = 132.2778
$ws.Range("N74").Value = -2802.4445
$ws.Range("H77").Value = 845.96295
$ws.Range("I77").Value = 741.7222
$ws.Range("J77").Value = 1054.4445
$ws.Range("K77").Value = 3708.611
$ws.Range("L77").Value = 5272.2225
$ws.Range("M77").Value = 659.3889999999997
$ws.Range("N77").Value = -14008.2225
$ws.Range("H110").Value = 167016880.0
$ws.Range("I110").Value = 167016880.0
$ws.Range("K110").Value = 167016880.0
$ws.Range("M110").Value = -167014835.0
$ws.Range("H111").Value = 26410.5
$ws.Range("J111").Value = 26410.5
$ws.Range("L111").Value = 26410.5
$ws.Range("N111").Value = -34590.5
$ws.Range("H112").Value = 9617.4
$ws.Range("J112").Value = 9617.4
$ws.Range("L112").Value = 9617.4
$ws.Range("N112").Value = -12571.4
$ws.Range("H132").Value = 3462.1702
$ws.Range("I132").Value = 3219.6155
$ws.Range("J132").Value = 4644.625
$ws.Range("K132").Value = 9658.8465
$ws.Range("L132").Value = 13933.875
$ws.Range("M132").Value = -7128.8465
$ws.Range("N132").Value = -18993.875

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 45775.4
$ws.Range("I86").Value = 56620.25
$ws.Range("J86").Value = 2396.0
$ws.Range("K86").Value = 56620.25
$ws.Range("L86").Value = 2396.0
$ws.Range("M86").Value = -55497.25
$ws.Range("N86").Value = -4642.0
$ws.Range("H89").Value = 45775.4
$ws.Range("I89").Value = 56620.25
$ws.Range("J89").Value = 2396.0
$ws.Range("K89").Value = 283101.25
$ws.Range("L89").Value = 11980.0
$ws.Range("M89").Value = -277485.25
$ws.Range("N89").Value = -23212.0
$ws.Range("H134").Value = 2181.8572
$ws.Range("I134").Value = 1922.1482
$ws.Range("J134").Value = 3058.375
$ws.Range("K134").Value = 5766.444600000001
$ws.Range("L134").Value = 9175.125
$ws.Range("M134").Value = -3231.444600000001
$ws.Range("N134").Value = -14245.125

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 108.57895
$ws.Range("I7").Value = 65.181816
$ws.Range("J7").Value = 168.25
$ws.Range("K7").Value = 65.181816
$ws.Range("L7").Value = 168.25
$ws.Range("M7").Value = 47.818184
$ws.Range("N7").Value = -394.25
$ws.Range("H32").Value = 22505.5
$ws.Range("J32").Value = 25011.0
$ws.Range("L32").Value = 25011.0
$ws.Range("N32").Value = -25643.0
$ws.Range("H37").Value = 41028.5
$ws.Range("J37").Value = 41028.5
$ws.Range("L37").Value = 41028.5
$ws.Range("N37").Value = -41242.5
$ws.Range("H58").Value = 2482.5557
$ws.Range("I58").Value = 2418.4285
$ws.Range("J58").Value = 2707.0
$ws.Range("K58").Value = 2418.4285
$ws.Range("L58").Value = 2707.0
$ws.Range("M58").Value = -2215.4285
$ws.Range("N58").Value = -3113.0
$ws.Range("H68").Value = 19958.4
$ws.Range("J68").Value = 19958.4
$ws.Range("L68").Value = 19958.4
$ws.Range("N68").Value = -21456.4
$ws.Range("H71").Value = 19958.4
$ws.Range("J71").Value = 19958.4
$ws.Range("L71").Value = 59875.2
$ws.Range("N71").Value = -67363.20000000001
$ws.Range("H74").Value = 31750.0
$ws.Range("J74").Value = 31750.0
$ws.Range("L74").Value = 31750.0
$ws.Range("N74").Value = -33498.0
$ws.Range("H77").Value = 31750.0
$ws.Range("J77").Value = 31750.0
$ws.Range("L77").Value = 95250.0
$ws.Range("N77").Value = -103986.0
$ws.Range("H88").Value = 22671.5
$ws.Range("I88").Value = 0.0
$ws.Range("K88").Value = 0.0
$ws.Range("M88").ClearContents()
$ws.Range("H91").Value = 22671.5
$ws.Range("I91").Value = 0.0
$ws.Range("K91").Value = 0.0
$ws.Range("M91").ClearContents()
$ws.Range("H124").Value = 22997.5
$ws.Range("J124").Value = 22997.5
$ws.Range("L124").Value = 22997.5
$ws.Range("N124").Value = -27907.5
$ws.Range("H132").Value = 3044.2424
$ws.Range("I132").Value = 3130.52
$ws.Range("J132").Value = 2774.625
$ws.Range("K132").Value = 9391.56
$ws.Range("L132").Value = 8323.875
$ws.Range("M132").Value = -6861.559999999999
$ws.Range("N132").Value = -13383.875
$ws.Range("H136").Value = 2482.5557
$ws.Range("I136").Value = 2418.4285
$ws.Range("J136").Value = 2707.0
$ws.Range("K136").Value = 7255.2855
$ws.Range("L136").Value = 8121.0
$ws.Range("M136").Value = -4705.2855
$ws.Range("N136").Value = -13221.0

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1260.3055
$ws.Range("I5").Value = 1307.25
$ws.Range("K5").Value = 3921.75
$ws.Range("M5").Value = -3809.75
$ws.Range("H7").Value = 667.3333
$ws.Range("J7").Value = 667.3333
$ws.Range("L7").Value = 2001.9999
$ws.Range("N7").Value = -2225.9999
$ws.Range("H12").Value = 41.4375
$ws.Range("J12").Value = 48.083332
$ws.Range("L12").Value = 144.249996
$ws.Range("N12").Value = -490.249996
$ws.Range("H34").Value = 1185.2
$ws.Range("I34").Value = 563.0
$ws.Range("J34").Value = 1600.0
$ws.Range("K34").Value = 1689.0
$ws.Range("L34").Value = 4800.0
$ws.Range("M34").Value = -1605.0
$ws.Range("N34").Value = -4968.0
$ws.Range("H58").Value = 2771.4285
$ws.Range("I58").Value = 2700.0
$ws.Range("J58").Value = 2950.0
$ws.Range("K58").Value = 8100.0
$ws.Range("L58").Value = 8850.0
$ws.Range("M58").Value = -7972.0
$ws.Range("N58").Value = -9106.0
$ws.Range("H92").Value = 629.0
$ws.Range("I92").Value = 400.0
$ws.Range("J92").Value = 800.75
$ws.Range("K92").Value = 1200.0
$ws.Range("L92").Value = 2402.25
$ws.Range("M92").Value = 48.0
$ws.Range("N92").Value = -4898.25
$ws.Range("H131").Value = 823.9899
$ws.Range("I131").Value = 383.0
$ws.Range("J131").Value = 837.7708
$ws.Range("K131").Value = 1149.0
$ws.Range("L131").Value = 2513.3124
$ws.Range("M131").Value = 3891.0
$ws.Range("N131").Value = -12593.3124
$ws.Range("H135").Value = 1260.3055
$ws.Range("I135").Value = 1307.25
$ws.Range("K135").Value = 11765.25
$ws.Range("M135").Value = -9230.25

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1464.9375
$ws.Range("I113").Value = 1242.6666
$ws.Range("J113").Value = 1598.3
$ws.Range("K113").Value = 1242.6666
$ws.Range("L113").Value = 1598.3
$ws.Range("M113").Value = 927.3334
$ws.Range("N113").Value = -5938.3

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 533281.56
$ws.Range("I46").Value = 320.0
$ws.Range("J46").Value = 1446930.0
$ws.Range("K46").Value = 320.0
$ws.Range("L46").Value = 1446930.0
$ws.Range("M46").Value = -132.0
$ws.Range("N46").Value = -1447306.0
$ws.Range("H127").Value = 31400.0
$ws.Range("J127").Value = 31400.0
$ws.Range("L127").Value = 31400.0
$ws.Range("N127").Value = -41320.0
$ws.Range("H128").Value = 18615.0
$ws.Range("J128").Value = 18615.0
$ws.Range("L128").Value = 18615.0
$ws.Range("N128").Value = -28575.0
$ws.Range("H132").Value = 3752.4075
$ws.Range("I132").Value = 4203.875
$ws.Range("K132").Value = 12611.625
$ws.Range("M132").Value = -10081.625

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 40000.0
$ws.Range("J31").Value = 40000.0
$ws.Range("L31").Value = 40000.0
$ws.Range("N31").Value = -40696.0
$ws.Range("H32").Value = 11762.917
$ws.Range("J32").Value = 12327.272
$ws.Range("L32").Value = 12327.272
$ws.Range("N32").Value = -12961.272
$ws.Range("H56").Value = 22047.5
$ws.Range("I56").Value = 3600.0
$ws.Range("J56").Value = 40495.0
$ws.Range("K56").Value = 3600.0
$ws.Range("L56").Value = 40495.0
$ws.Range("M56").Value = -2886.0
$ws.Range("N56").Value = -41923.0
$ws.Range("H62").Value = 5130491.5
$ws.Range("I62").Value = 12822345.0
$ws.Range("J62").Value = 2588.889
$ws.Range("K62").Value = 12822345.0
$ws.Range("L62").Value = 2588.889
$ws.Range("M62").Value = -12821721.0
$ws.Range("N62").Value = -3836.889
$ws.Range("H65").Value = 5130491.5
$ws.Range("I65").Value = 12822345.0
$ws.Range("J65").Value = 2588.889
$ws.Range("K65").Value = 64111725.0
$ws.Range("L65").Value = 12944.445
$ws.Range("M65").Value = -64108605.0
$ws.Range("N65").Value = -19184.445
$ws.Range("H122").Value = 3220.7693
$ws.Range("I122").Value = 2180.0
$ws.Range("J122").Value = 3871.25
$ws.Range("K122").Value = 6540.0
$ws.Range("L122").Value = 11613.75
$ws.Range("M122").Value = -4090.0
$ws.Range("N122").Value = -16513.75
$ws.Range("H124").Value = 29417.555
$ws.Range("J124").Value = 29417.555
$ws.Range("L124").Value = 29417.555
$ws.Range("N124").Value = -39237.555
